$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-12-04 Thursday" "2025-12-05 Friday"

Replace-Text "880÷6=146, 4" "105÷7=15, 0"
Replace-Text "847÷8=105, 7" "473÷8=59, 1"
Replace-Text "219÷5=43, 4" "424÷6=70, 4"
Replace-Text "301÷3=100, 1" "963÷8=120, 3"
Replace-Text "283÷9=31, 4" "576÷2=288, 0"

Replace-Text "664÷6=110, 4" "237÷2=118, 1"
Replace-Text "499÷2=249, 1" "711÷7=101, 4"
Replace-Text "917÷6=152, 5" "382÷2=191, 0"
Replace-Text "830÷6=138, 2" "462÷5=92, 2"
Replace-Text "802÷2=401, 0" "558÷8=69, 6"

Replace-Text "418÷5=83, 3" "841÷4=210, 1"
Replace-Text "376÷7=53, 5" "574÷8=71, 6"
Replace-Text "280÷2=140, 0" "641÷8=80, 1"
Replace-Text "906÷5=181, 1" "270÷7=38, 4"
Replace-Text "904÷4=226, 0" "207÷4=51, 3"

Replace-Text "308÷2=154, 0" "474÷8=59, 2"
Replace-Text "737÷9=81, 8" "296÷5=59, 1"
Replace-Text "897÷6=149, 3" "322÷4=80, 2"
Replace-Text "883÷9=98, 1" "692÷5=138, 2"
Replace-Text "313÷3=104, 1" "425÷7=60, 5"

Replace-Text "991÷5=198, 1" "610÷7=87, 1"
Replace-Text "441÷6=73, 3" "106÷4=26, 2"
Replace-Text "747÷9=83, 0" "645÷3=215, 0"
Replace-Text "237÷8=29, 5" "965÷9=107, 2"
Replace-Text "734÷5=146, 4" "807÷9=89, 6"

Write-Host "Done"
